# Apply the UC1_UC1 error-log edit: rename the user on every row, renumber
# the screenshot paths under the new capture folder, and rewrite several
# explanation strings. Row 5/7 additionally swap type/error content to
# describe a new Windows Update error (0x80240fff) replacing the old one
# (0x80244007), which moves from row 7 to row 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: user_name, every data row 2..16 gets renamed ---
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 3).Value = "Naoto Ikeda"
}

# --- Column J: capimg path, renumbered per-row under the new folder ---
$capimg = @{
    2  = "bdot20240415_141954/1.png"
    3  = "bdot20240415_141954/2.png"
    4  = "bdot20240415_141954/3.png"
    5  = "bdot20240415_141954/4.png"
    6  = "bdot20240415_141954/5.png"
    7  = "bdot20240415_141954/5.png"
    8  = "bdot20240415_141954/6.png"
    9  = "bdot20240415_141954/7.png"
    10 = "bdot20240415_141954/8.png"
    11 = "bdot20240415_141954/9.png"
    12 = "bdot20240415_141954/10.png"
    13 = "bdot20240415_141954/1.png"
    14 = "bdot20240415_141954/2.png"
    15 = "bdot20240415_141954/3.png"
    16 = "bdot20240415_141954/11.png"
}
foreach ($r in $capimg.Keys) {
    $ws.Cells.Item($r, 10).Value = $capimg[$r]
}

# --- Column K: explanation text, rewritten per-row ---
$explanation = @{
    2  = "「スタート」ボタンをクリックする"
    3  = "メニューから「設定」アイコンをクリックする"
    4  = "左側のメニューからWindows Updateをクリックし、Windows Update画面に移動する"
    5  = "0x80240fff エラー"
    6  = "デスクトップ画面の左下にある「スタート」ボタンを右クリックする"
    7  = "メニューからターミナル(管理者)をクリックする"
    8  = "ユーザーアカウント制御と表示されているウィンドウが開いたことを確認する"
    9  = "PowerShellウィンドウに start-transcript と入力し、[Enter]キーを押す"
    10 = "wuauclt.exe /resetauthorization /detectnow と入力し、[Enter]キーを押す"
    11 = "netsh winhttp show proxy と入力し、[Enter]キーを押す"
    12 = "netsh winhttp reset proxy と入力し、[Enter]キーを押す"
    13 = "「スタート」ボタンをクリックする"
    14 = "メニューから「設定」アイコンをクリックする"
    15 = "左側のメニューからWindows Updateをクリックし、Windows Update画面に移動する"
    16 = "「更新プログラムのチェック」ボタンをクリックする"
}
foreach ($r in $explanation.Keys) {
    $ws.Cells.Item($r, 11).Value = $explanation[$r]
}

# --- Row 5 becomes "error" type, gains error_type/error_content ---
$ws.Cells.Item(5, 2).Value = "error"
$ws.Cells.Item(5, 12).Value = "Error W"
$ws.Cells.Item(5, 13).Value = " エラーの Windows"

# --- Row 7 becomes "operation" type, loses error_type/error_content ---
$ws.Cells.Item(7, 2).Value = "operation"
$ws.Cells.Item(7, 12).Value = ""
$ws.Cells.Item(7, 13).Value = ""
